# Colando o meu código teste que funcionou no Main
# Append two new FIPE rows for Agrale "MARRUÁ AM 100 2.8  CS TDI Diesel"
# (2015 and 2014 Diesel) right after the existing data (rows 2-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - 2015 Diesel
$ws.Range("A21").Value = "Agrale"
$ws.Range("B21").Value = "MARRUÁ AM 100 2.8  CS TDI Diesel"
$ws.Range("C21").Value = "2015 Diesel"
$ws.Range("D21").Value = "060003-2"
$ws.Range("E21").Value = "'" + " 108542.00"
$ws.Range("E21").Style = "Normal"

# Row 22 - 2014 Diesel
$ws.Range("A22").Value = "Agrale"
$ws.Range("B22").Value = "MARRUÁ AM 100 2.8  CS TDI Diesel"
$ws.Range("C22").Value = "2014 Diesel"
$ws.Range("D22").Value = "060003-2"
$ws.Range("E22").Value = "'" + " 102457.00"
$ws.Range("E22").Style = "Normal"
